# Global Irminger - fixed deployment dates
# Fixed deployment dates based on cruise reports and WHOI documentation

$wb = $excel.ActiveWorkbook

$wsMoorings = $wb.Worksheets.Item("Moorings")

# Recover Date (G2) was missing - fill in from cruise report (2015-11-05)
$wsMoorings.Range("G2").Value = 42313

# Notes (L2) - the glider was lost, so note it
$wsMoorings.Range("L2").Value = "glider lost"

# The Moorings sheet is now the one being worked on / reviewed
$wsMoorings.Activate() | Out-Null
$wsMoorings.Range("F11").Select() | Out-Null
